# Colocando header nos gráficos
# Adds a header label to column A (row 1) of each data sheet feeding the
# charts, fixes accented Portuguese text, and removes the now-unused
# "Teto" row from the emissions sheet / refreshes the cost sheet numbers.

$wb = $excel.ActiveWorkbook

# xlPasteFormats=-4122, xlPasteValues=-4163
$xlPasteFormats = -4122
$xlPasteValues  = -4163

function Set-HeaderCell {
    param($ws, [string]$text)

    # Grab the existing header style (already bold/centered/bordered on B1)
    # and stamp it onto A1 without creating a brand-new style entry.
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("A1").Value = $text
}

function Copy-TextValue {
    # Copies the VALUE ONLY of $sourceRange into $destCell so a numeric-
    # looking label (e.g. the year "2015") lands as genuine text instead
    # of falling prey to Excel's auto-numeric-coercion on assignment --
    # while leaving $destCell's existing style (and cellXf index) intact.
    param($sourceRange, $destCell)

    $sourceRange.Copy() | Out-Null
    $destCell.PasteSpecial($xlPasteValues) | Out-Null
}

# ---------------------------------------------------------------------
# Sheets 1-4: "Fonte/Tecnologia" header + accent fixes, row labels lose
# their bold/bordered style (s="1" removed -> default style).
# ---------------------------------------------------------------------
$sourceSheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

$labels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

foreach ($sheetName in $sourceSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    Set-HeaderCell $ws "Fonte/Tecnologia"

    foreach ($row in 2..12) {
        $cell = $ws.Cells.Item($row, 1)
        $cell.Value = $labels[$row]
        $cell.ClearFormats() | Out-Null
    }
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)" -> "Período" header, accent fixes,
# and removal of the "Teto" row (row 4).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

Set-HeaderCell $ws5 "Período"

$ws5.Cells.Item(2, 1).Value = "P.Médio"
$ws5.Cells.Item(2, 1).ClearFormats() | Out-Null

$ws5.Cells.Item(3, 1).Value = "P.Crítico"
$ws5.Cells.Item(3, 1).ClearFormats() | Out-Null

$ws5.Rows.Item(4).Delete() | Out-Null

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)" -> "Tipo Expansão" header,
# B1 relabelled "2015", accent fixes and refreshed values.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

Set-HeaderCell $ws6 "Tipo Expansão"

# Source the literal text "2015" from another sheet's year-header cell
# (already a text value there) so it lands here as text, not a number.
$yearTextSource = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)").Range("B1")
Copy-TextValue $yearTextSource $ws6.Range("B1")

$ws6.Cells.Item(2, 1).Value = "Expansão Centralizada"
$ws6.Cells.Item(2, 1).ClearFormats() | Out-Null
$ws6.Cells.Item(2, 2).Value = 609

$ws6.Cells.Item(3, 1).Value = "Expansão por GD"
$ws6.Cells.Item(3, 1).ClearFormats() | Out-Null
$ws6.Cells.Item(3, 2).Value = 99

Write-Output "done"
